# Daily attendance processing - 2026-01-06 15:38:15
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (G) wherever both are listed as "System, <email>".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns("G")
$null = $col.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
